$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header for column C from "4 Owned" to "# Owned"
$ws.Range("C1").Value = "# Owned"

# Update column C values from text "Yes" to numeric "owned" counts
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 4

# Move the active selection to C3
$ws.Range("C3").Select()
